# Update the 预住院总结报表 workbook: a new batch of records was folded
# into every pre-aggregated summary sheet, raising the year total from
# 1352 to 1610 (a +258 delta that shows up consistently across all
# twelve views of the same underlying data).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 预住院年流量 (year total) ---------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = 1610

# --- Sheet 2: 月流量 (by month, rows) -----------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A12").Value = 11
$ws2.Range("B12").Value = 258

# --- Sheet 3: 周流量 (by week number, rows) -----------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B43").Value = 56
$ws3.Range("A44").Value = 45
$ws3.Range("B44").Value = 59
$ws3.Range("A45").Value = 46
$ws3.Range("B45").Value = 75
$ws3.Range("A46").Value = 47
$ws3.Range("B46").Value = 64
$ws3.Range("A47").Value = 48
$ws3.Range("B47").Value = 52

# --- Sheet 4: 周内流量 (by weekday, rows) -------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = 356
$ws4.Range("B3").Value = 267
$ws4.Range("B4").Value = 303
$ws4.Range("B5").Value = 235
$ws4.Range("B6").Value = 252
$ws4.Range("B7").Value = 140
$ws4.Range("B8").Value = 57

# --- Sheet 5: 核算年 (year total, accounting view) ----------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = 1610

# --- Sheet 6: 核算月 (by month, accounting view, rows) ------------------
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("B12").Value = 280
$ws6.Range("A13").Value = 12
$ws6.Range("B13").Value = 32

# --- Sheet 7: 核算年核算月流量 (by month, accounting view, cols) --------
$ws7 = $wb.Worksheets.Item(7)
$ws7.Cells.Item(2, 12).Value = 280
$ws7.Cells.Item(1, 13).Value = 12
$ws7.Cells.Item(2, 13).Value = 32

# --- Sheet 8: 年周期月度流量 (by month, cols) ----------------------------
$ws8 = $wb.Worksheets.Item(8)
$ws8.Cells.Item(1, 12).Value = 11
$ws8.Cells.Item(2, 12).Value = 258

# --- Sheet 9: 年周期周度流量 (by week number, cols) ----------------------
$ws9 = $wb.Worksheets.Item(9)
$ws9.Cells.Item(2, 43).Value = 56
$ws9.Cells.Item(1, 44).Value = 45
$ws9.Cells.Item(2, 44).Value = 59
$ws9.Cells.Item(1, 45).Value = 46
$ws9.Cells.Item(2, 45).Value = 75
$ws9.Cells.Item(1, 46).Value = 47
$ws9.Cells.Item(2, 46).Value = 64
$ws9.Cells.Item(1, 47).Value = 48
$ws9.Cells.Item(2, 47).Value = 52

# --- Sheet 10: 年周期月中流量 (by weekday-in-month, cols) -----------------
$ws10 = $wb.Worksheets.Item(10)
$newVals10 = @(31,24,21,50,36,50,45,52,48,37,79,64,61,52,65,49,42,64,58,66,56,83,55,62,62,71,45,56,52,42,32)
$col = 2
foreach ($v in $newVals10) {
  $ws10.Cells.Item(2, $col).Value = $v
  $col = $col + 1
}

# --- Sheet 11: 年周期周中流量 (by weekday, cols) ---------------------------
$ws11 = $wb.Worksheets.Item(11)
$ws11.Range("B2").Value = 356
$ws11.Range("C2").Value = 267
$ws11.Range("D2").Value = 303
$ws11.Range("E2").Value = 235
$ws11.Range("F2").Value = 252
$ws11.Range("G2").Value = 140
$ws11.Range("H2").Value = 57

# --- Sheet 12: 月周期周中流量 (by week-in-month x weekday, rows) -----------
$ws12 = $wb.Worksheets.Item(12)
$ws12.Range("A12").Value = 11
$ws12.Range("B12").Value = 75
$ws12.Range("C12").Value = 36
$ws12.Range("D12").Value = 37
$ws12.Range("E12").Value = 35
$ws12.Range("F12").Value = 44
$ws12.Range("G12").Value = 22
$ws12.Range("H12").Value = 9
